# Weekly update: insert a new Perejil (Mercado Mayorista Lo Valledor de Santiago)
# record ahead of the existing rows, pushing the rest of the table down by one row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 309 - everything currently at/after row 309
# (through the old last row 339) shifts down to 310..340.
$ws.Rows.Item(309).Insert()

# Populate the newly inserted row 309 with the new weekly observation.
$ws.Cells.Item(309, 1).Value = 6
$ws.Cells.Item(309, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(309, 3).Value = "Metropolitana"
$ws.Cells.Item(309, 4).Value = 44449
$ws.Cells.Item(309, 5).Value = 13
$ws.Cells.Item(309, 6).Value = 100112044
$ws.Cells.Item(309, 7).Value = "Perejil"
$ws.Cells.Item(309, 8).Value = "Sin especificar"
$ws.Cells.Item(309, 9).Value = "Primera"
$ws.Cells.Item(309, 10).Value = 230
$ws.Cells.Item(309, 11).Value = 8000
$ws.Cells.Item(309, 12).Value = 9000
$ws.Cells.Item(309, 13).Value = 8391
$ws.Cells.Item(309, 14).Value = "`$/docena de atados"
$ws.Cells.Item(309, 15).Value = "Región Metropolitana"
$ws.Cells.Item(309, 16).Value = 2797
$ws.Cells.Item(309, 17).Value = 3
$ws.Cells.Item(309, 18).Value = "Hortaliza"
